$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-label biological network categories: "Neural", "Protein", and
#    "Genetic" are being consolidated into a single "Biological" label
#    (column M, "Label").
$usedRows = $ws.UsedRange.Rows.Count
for ($i = 1; $i -le $usedRows; $i++) {
    $cell = $ws.Cells.Item($i, 13)
    $val = $cell.Value2
    if ($val -eq "Neural" -or $val -eq "Protein" -or $val -eq "Genetic") {
        $cell.Value2 = "Biological"
    }
}

# 2) Remove the Celegans_Multiplex_Genetic (rows 71-72) and
#    Drosophila_Multiplex_Genetic (rows 73-75) data rows entirely - the
#    three-layer C++ experiment replaces them, so the old two datasets'
#    rows are dropped and everything below shifts up by 5 rows.
$ws.Range("A71:M75").EntireRow.Delete()

# 3) Update the view: scroll the frozen pane down and move the active
#    selection onto the newly-relevant rows.
$win = $excel.Application.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("A71:XFD73").Select()
